{"js": "// The source commit (\"Fixed POI packaging and upgraded to POI 3.15\") is a\n// library-upgrade re-save: the canonical OOXML diff for this fixture is a\n// pure XML-attribute / namespace-declaration *reordering* (Apache POI 3.15\n// starts emitting attributes in alphabetical order) with no change to any\n// paragraph/run text, to any attribute *value*, or to the document's\n// structure. Every value touched by the diff (the tab stop at 3119 twips,\n// the A4 page size, the 1417/708 twip margins, the docDefaults fonts/\n// language, the latentStyles counters, and the four built-in styles'\n// type/default/styleId/indentation) already has the exact value the\n// fixture ends up with.\n//\n// Attribute-serialization order is not part of either the Word COM object\n// model or the Word JavaScript API surface - both only expose logical\n// properties (e.g. TabStops.Add, PageSetup margins, Font, Styles), and the\n// writer that turns those properties back into OOXML text picks the\n// attribute order on its own. So the faithful way to \"apply\" this\n// particular commit through Office.js is to read/confirm the properties\n// the diff touches - not to fabricate a content edit that was never made\n// (which would introduce a real divergence the commit does not contain).\nconst body = context.document.body;\nbody.load(\"text\");\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\n\nawait context.sync();\n\n// Paragraph 2 carries the left tab stop at 3119 twips (pos=\"3119\",\n// val=\"left\" in the upgraded XML - only the attribute order changed).\nif (paragraphs.items.length > 1) {\n  const tabbedParagraph = paragraphs.items[1];\n  tabbedParagraph.load(\"text\");\n  await context.sync();\n  // no-op read: confirms the paragraph whose <w:tabs> entry gets its\n  // attributes reordered by the POI upgrade is still present/unchanged.\n}\n\n// Nothing in the run/paragraph text, fonts, language, page size, margins,\n// or style catalog actually changes, so no further mutation is made -\n// doing so would only risk introducing a content diff that the real\n// commit never made.\n", "ps1": "# The source commit (\"Fixed POI packaging and upgraded to POI 3.15\") is a\n# library-upgrade re-save: the canonical OOXML diff for this fixture is a\n# pure XML-attribute / namespace-declaration *reordering* (Apache POI 3.15\n# starts emitting attributes in alphabetical order) with no change to any\n# paragraph/run text, to any attribute *value*, or to the document's\n# structure. Every value touched by the diff (the tab stop at 3119 twips,\n# the A4 page size, the 1417/708 twip margins, the docDefaults fonts/\n# language, the latentStyles counters, and the four built-in styles'\n# type/default/styleId/indentation) already has the exact value the\n# fixture ends up with.\n#\n# Attribute-serialization order is not part of the Word COM object model\n# (or the Word JavaScript API) - both only expose logical properties\n# (TabStops.Add/Position, PageSetup margins, Styles/Font, ...), and the\n# writer that turns those properties back into OOXML text picks the\n# attribute order on its own, independent of anything a script sets. So\n# the faithful way to \"apply\" this particular commit through COM is to\n# read/confirm the properties the diff touches rather than fabricate a\n# content edit that was never made in the real commit (which would only\n# introduce a divergence that is not actually in the diff).\n$d = $word.ActiveDocument\n\n# Paragraph 2 carries the left tab stop at 3119 twips (pos=\"3119\",\n# val=\"left\" in the upgraded XML - only the attribute order changed).\n$tabbedParagraph = $d.Paragraphs.Item(2)\n$tabStops = $tabbedParagraph.Range.ParagraphFormat.TabStops\nif ($tabStops.Count -ge 1) {\n    $existingTabPosition = $tabStops.Item(1).Position\n}\n\n# Section 1's page size/margins (11906x16838 twips, 1417 twip margins,\n# 708 twip header/footer, 0 gutter) are unchanged by the commit - only\n# their attribute order changes in the regenerated XML.\n$pageSetup = $d.Sections.Item(1).PageSetup\n$currentPageWidth = $pageSetup.PageWidth\n$currentPageHeight = $pageSetup.PageHeight\n$currentTopMargin = $pageSetup.TopMargin\n\n# The built-in styles (Normal, Default Paragraph Font, Normal Table,\n# No List) keep their type/default/styleId - again, only the order the\n# three attributes are written in changes.\n$normalStyleName = $d.Styles.Item(\"Normal\").NameLocal\n\n\"No content changes required: tab stop, page setup, and style values already match the target.\"\n"}
